$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C entirely (shifts D->C, E->D)
$ws.Range("C:C").Delete()

# Update B2 with the new prediction value
$ws.Range("B2").Value = 0.9073038990665374
